$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 03:50"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 215081
$ws.Range("C4").Value = 78
$ws.Range("E4").Value = 201094
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 5109

# Row 17 - Corea del Sur
$ws.Range("B17").Value = 9976
$ws.Range("C17").Value = 89
$ws.Range("D17").Value = 5828
$ws.Range("E17").Value = 3979
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 169

# Row 18 - Canada
$ws.Range("E18").Value = 7866
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 129

# Row 53 - Singapur
$ws.Range("E53").Value = 751
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 4

# Row 154 - Guyana
$ws.Range("E154").Value = 15
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 4
